# Updates the cryptos list (rows 2-51) to the new values scraped on
# Wed Aug 16 09:28:56 UTC 2023, matching the commit's XML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.241.39"
$ws.Range("E2").Value = "  -0.59%  "
# Row 3
$ws.Range("D3").Value = "1.828.82"
$ws.Range("E3").Value = "  -0.76%  "
# Row 4
$ws.Range("E4").Value = "  +0.53%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.76"
$ws.Range("E5").Value = "  -1.96%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6027"
$ws.Range("E6").Value = "  -4.01%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.41%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07024"
$ws.Range("E8").Value = "  -5.65%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2788"
$ws.Range("E9").Value = "  -3.78%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.36"
$ws.Range("E10").Value = "  -6.36%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07670"
$ws.Range("E11").Value = "  -0.62%  "
# Row 12
$ws.Range("D12").Value = "1.822.60"
$ws.Range("E12").Value = "  -1.01%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.794"
$ws.Range("E13").Value = "  -3.63%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.000009884"
$ws.Range("E14").Value = "  -3.75%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6234"
$ws.Range("E15").Value = "  -7.86%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "79.03"
$ws.Range("E16").Value = "  -3.46%  "
# Row 17
$ws.Range("D17").Value = "29.238.20"
$ws.Range("E17").Value = "  -0.77%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.827"
$ws.Range("E18").Value = "  -6.67%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "224.24"
$ws.Range("E19").Value = "  -3.80%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.37%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.67"
$ws.Range("E21").Value = "  -5.31%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.993"
$ws.Range("E22").Value = "  -4.55%  "
# Row 23
$ws.Range("E23").Value = "  +0.46%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "156.62"
$ws.Range("E24").Value = "  -1.00%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.989"
$ws.Range("E25").Value = "  -5.97%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1292"
$ws.Range("E26").Value = "  -4.42%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.52"
$ws.Range("E27").Value = "  -4.85%  "
# Row 28
$ws.Range("E28").Value = "  +0.78%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06182"
$ws.Range("E29").Value = "  -13.22%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.441"
$ws.Range("E30").Value = "  -2.86%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.819"
$ws.Range("E31").Value = "  -5.47%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.789"
$ws.Range("E32").Value = "  -6.38%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.118"
$ws.Range("E33").Value = "  -1.95%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.741"
$ws.Range("E34").Value = "  -4.35%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6428"
$ws.Range("E35").Value = "  -7.62%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.544"
$ws.Range("E36").Value = "  -1.24%  "
# Row 37
$ws.Range("D37").Value = "1.219.54"
$ws.Range("E37").Value = "  -1.37%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.739"
$ws.Range("E38").Value = "  -2.80%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.547"
$ws.Range("E39").Value = "  -6.00%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01731"
$ws.Range("E40").Value = "  -5.99%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8977"
$ws.Range("E41").Value = "  -6.24%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("E42").Value = "  +0.40%  "
# Row 43
$ws.Range("D43").Value = "1.987.48"
$ws.Range("E43").Value = "  -0.90%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.34"
$ws.Range("E44").Value = "  -0.64%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.40"
$ws.Range("E45").Value = "  -4.74%  "
# Row 46
$ws.Range("E46").Value = "  -3.05%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.496"
$ws.Range("E47").Value = "  -5.01%  "
# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4563"
$ws.Range("E48").Value = "  -0.43%  "
# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.572"
$ws.Range("E49").Value = "  -9.13%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05507"
$ws.Range("E50").Value = "  -2.65%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.407"
$ws.Range("E51").Value = "  -8.06%  "
